$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 4993.8887
$ws.Range("I8").Value = 2368.125
$ws.Range("K8").Value = 7104.375
$ws.Range("M8").Value = -6965.375
$ws.Range("H112").Value = 2075.7058
$ws.Range("I112").Value = 549.5
$ws.Range("J112").Value = 2545.3076
$ws.Range("K112").Value = 1648.5
$ws.Range("L112").Value = 7635.9228
$ws.Range("M112").Value = -540.5
$ws.Range("N112").Value = -9851.9228
$ws.Range("H129").Value = 1169.4717
$ws.Range("I129").Value = 367.5
$ws.Range("J129").Value = 1234.9387
$ws.Range("K129").Value = 1102.5
$ws.Range("L129").Value = 3704.8161
$ws.Range("M129").Value = 3897.5
$ws.Range("N129").Value = -13704.8161
$ws.Range("H132").Value = 1798.3158
$ws.Range("I132").Value = 1457
$ws.Range("J132").Value = 3888.875
$ws.Range("K132").Value = 4371
$ws.Range("L132").Value = 11666.625
$ws.Range("M132").Value = -1841
$ws.Range("N132").Value = -16726.625
$ws.Range("H137").Value = 1602.6666
$ws.Range("I137").Value = 1511.3846
$ws.Range("J137").Value = 1840
$ws.Range("K137").Value = 4534.1538
$ws.Range("L137").Value = 5520
$ws.Range("M137").Value = -1984.1538
$ws.Range("N137").Value = -10620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16976.645
$ws.Range("I32").Value = 18859
$ws.Range("J32").Value = 4741.3335
$ws.Range("K32").Value = 18859
$ws.Range("L32").Value = 4741.3335
$ws.Range("M32").Value = -18572
$ws.Range("N32").Value = -5315.3335
$ws.Range("H61").Value = 2246.318
$ws.Range("I61").Value = 2095.95
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 2095.95
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -1883.95
$ws.Range("N61").Value = -4174
$ws.Range("H74").Value = 1563.091
$ws.Range("I74").Value = 900
$ws.Range("J74").Value = 1942
$ws.Range("K74").Value = 900
$ws.Range("L74").Value = 1942
$ws.Range("M74").Value = -26
$ws.Range("N74").Value = -3690
$ws.Range("H77").Value = 1563.091
$ws.Range("I77").Value = 900
$ws.Range("J77").Value = 1942
$ws.Range("K77").Value = 4500
$ws.Range("L77").Value = 9710
$ws.Range("M77").Value = -132
$ws.Range("N77").Value = -18446
$ws.Range("H122").Value = 1353.1154
$ws.Range("I122").Value = 1223.24
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 3669.72
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -1219.72
$ws.Range("N122").Value = -18700
$ws.Range("H132").Value = 9271.117
$ws.Range("I132").Value = 12234.637
$ws.Range("J132").Value = 3838
$ws.Range("K132").Value = 36703.911
$ws.Range("L132").Value = 11514
$ws.Range("M132").Value = -34173.911
$ws.Range("N132").Value = -16574
$ws.Range("H136").Value = 2246.318
$ws.Range("I136").Value = 2095.95
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 6287.849999999999
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -3737.849999999999
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 63899.75
$ws.Range("I20").Value = 68059.734
$ws.Range("K20").Value = 68059.734
$ws.Range("M20").Value = -67812.734
$ws.Range("H22").Value = 397.07693
$ws.Range("I22").Value = 379.5
$ws.Range("J22").Value = 425.2
$ws.Range("K22").Value = 379.5
$ws.Range("L22").Value = 425.2
$ws.Range("M22").Value = -206.5
$ws.Range("N22").Value = -771.2
$ws.Range("H86").Value = 145278.14
$ws.Range("I86").Value = 2788.4
$ws.Range("J86").Value = 501502.5
$ws.Range("K86").Value = 2788.4
$ws.Range("L86").Value = 501502.5
$ws.Range("M86").Value = -1665.4
$ws.Range("N86").Value = -503748.5
$ws.Range("H89").Value = 145278.14
$ws.Range("I89").Value = 2788.4
$ws.Range("J89").Value = 501502.5
$ws.Range("K89").Value = 13942
$ws.Range("L89").Value = 2507512.5
$ws.Range("M89").Value = -8326
$ws.Range("N89").Value = -2518744.5
$ws.Range("H105").Value = 3031.92
$ws.Range("I105").Value = 2957.1765
$ws.Range("J105").Value = 3190.75
$ws.Range("K105").Value = 2957.1765
$ws.Range("L105").Value = 3190.75
$ws.Range("M105").Value = -1210.1765
$ws.Range("N105").Value = -6684.75
$ws.Range("H134").Value = 3242.5833
$ws.Range("I134").Value = 3391.1
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 10173.3
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -7638.299999999999
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2907.9656
$ws.Range("I31").Value = 2100.8235
$ws.Range("J31").Value = 4051.4167
$ws.Range("K31").Value = 2100.8235
$ws.Range("L31").Value = 4051.4167
$ws.Range("M31").Value = -1805.8235
$ws.Range("N31").Value = -4641.4167
$ws.Range("H34").Value = 2907.9656
$ws.Range("I34").Value = 2100.8235
$ws.Range("J34").Value = 4051.4167
$ws.Range("K34").Value = 2100.8235
$ws.Range("L34").Value = 4051.4167
$ws.Range("M34").Value = -1898.8235
$ws.Range("N34").Value = -4455.4167
$ws.Range("H122").Value = 1844.5714
$ws.Range("I122").Value = 1927.1765
$ws.Range("J122").Value = 1493.5
$ws.Range("K122").Value = 5781.529500000001
$ws.Range("L122").Value = 4480.5
$ws.Range("M122").Value = -3331.529500000001
$ws.Range("N122").Value = -9380.5
$ws.Range("H132").Value = 589831.8
$ws.Range("I132").Value = 677205.9399999999
$ws.Range("K132").Value = 2031617.82
$ws.Range("M132").Value = -2029087.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 297.25
$ws.Range("I10").Value = 69.666664
$ws.Range("J10").Value = 980
$ws.Range("K10").Value = 208.999992
$ws.Range("L10").Value = 2940
$ws.Range("M10").Value = -69.99999199999999
$ws.Range("N10").Value = -3218
$ws.Range("H14").Value = 41
$ws.Range("I14").Value = 41
$ws.Range("K14").Value = 123
$ws.Range("M14").Value = 50
$ws.Range("H62").Value = 1000
$ws.Range("J62").Value = 1000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4372
$ws.Range("H63").Value = 236445.89
$ws.Range("H64").Value = 3560.5454
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3560.5454
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10681.6362
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -11221.6362
$ws.Range("H65").Value = 1000
$ws.Range("J65").Value = 1000
$ws.Range("L65").Value = 9000
$ws.Range("N65").Value = -15864
$ws.Range("H66").Value = 236445.89
$ws.Range("H67").Value = 3560.5454
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3560.5454
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10681.6362
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -12553.6362
$ws.Range("H68").Value = 594.3333
$ws.Range("J68").Value = 503
$ws.Range("L68").Value = 1509
$ws.Range("N68").Value = -3131
$ws.Range("H71").Value = 594.3333
$ws.Range("J71").Value = 503
$ws.Range("L71").Value = 4527
$ws.Range("N71").Value = -12639
$ws.Range("H122").Value = 869
$ws.Range("J122").Value = 1999
$ws.Range("L122").Value = 17991
$ws.Range("N122").Value = -22891

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1147
$ws.Range("I16").Value = 1173.8823
$ws.Range("J16").Value = 690
$ws.Range("K16").Value = 1173.8823
$ws.Range("L16").Value = 690
$ws.Range("M16").Value = -1003.8823
$ws.Range("N16").Value = -1030
$ws.Range("H22").Value = 542.625
$ws.Range("I22").Value = 505.85715
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 505.85715
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -210.85715
$ws.Range("N22").Value = -1390
$ws.Range("H27").Value = 542.625
$ws.Range("I27").Value = 505.85715
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 505.85715
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -398.85715
$ws.Range("N27").Value = -1014
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350
$ws.Range("H132").Value = 4395.793
$ws.Range("I132").Value = 4210.769
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 12632.307
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -10102.307
$ws.Range("N132").Value = -23058.0005
$ws.Range("H136").Value = 2780375.8
$ws.Range("I136").Value = 4466225
$ws.Range("J136").Value = 3682.353
$ws.Range("K136").Value = 13398675
$ws.Range("L136").Value = 11047.059
$ws.Range("M136").Value = -13396125
$ws.Range("N136").Value = -16147.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50406
